$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 169, pushing the existing row 169 (and everything
# below it, down through row 182) down to row 170 (through row 183).
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new "Albahaca" record.
$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value = "Los Lagos"
$ws.Cells.Item(169, 4).Value = 45013
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = 100112052
$ws.Cells.Item(169, 7).Value = "Albahaca"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 120
$ws.Cells.Item(169, 11).Value = 6500
$ws.Cells.Item(169, 12).Value = 7000
$ws.Cells.Item(169, 13).Value = 6750
$ws.Cells.Item(169, 14).Value = "$/docena de matas"
$ws.Cells.Item(169, 15).Value = "Región Metropolitana"
$ws.Cells.Item(169, 16).Value = 1125
$ws.Cells.Item(169, 17).Value = 6
$ws.Cells.Item(169, 18).Value = "Hortaliza"
